# Replace the per-scenario cost-breakdown tables (row 1 headers + row 2 values)
# on each year sheet with the latest results received from the server.
# Header layout changes too: "gt"/"dgt" columns are dropped and new "gb"
# (after "eb") and "btes" (after "ttes") columns are introduced.
$wb = $excel.ActiveWorkbook

$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

$ws = $wb.Worksheets.Item("2025")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("A2").Value = 3906.399109145206
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 48353.76274462014
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 9433.134471502228
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 2534.277928792104
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 2367.37219622158
$ws.Range("O2").Value = 1995.762462679798

$ws = $wb.Worksheets.Item("2030")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("A2").Value = 6991.052031681918
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 197913.7502057619
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16452.51445364119
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 8194.52068131253
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 7544.284082242735
$ws.Range("O2").Value = 6258.021366348374

$ws = $wb.Worksheets.Item("2035")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("A2").Value = 31236.29455387744
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 12889.13275549508
$ws.Range("O2").Value = 9263.90107805635

$ws = $wb.Worksheets.Item("2040")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("A2").Value = 31236.29455387744
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 14046.98250793826
$ws.Range("O2").Value = 9263.90107805635

$ws = $wb.Worksheets.Item("2045")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("A2").Value = 38906.8534480406
$ws.Range("B2").Value = 193.0947398408091
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 16880.98779587899
$ws.Range("O2").Value = 10096.4577740545

$ws = $wb.Worksheets.Item("2050")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("A2").Value = 38906.8534480406
$ws.Range("B2").Value = 193.0947398408091
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 16880.98779587899
$ws.Range("O2").Value = 10096.4577740545
